$d = $word.ActiveDocument

# 1. "a été créé via ce moyen." -> "a été créé par ce moyen."
$d.Content.Find.Execute("a été créé via ce moyen.", $false, $false, $false, $false, $false, $true, 1, $false, "a été créé par ce moyen.", 2) | Out-Null

# 2. "savoir à où en est le processus" -> "savoir où en est le processus"
$d.Content.Find.Execute("savoir à où en est le processus", $false, $false, $false, $false, $false, $true, 1, $false, "savoir où en est le processus", 2) | Out-Null

# 3. "il est ajouté à un calendrier. Celui-ci ... fusionner." -> "...il est ajouté au calendrier final. Celui-ci ... fusionner. Le résultat est disponible pour l'exportation."
$d.Content.Find.Execute("il est ajouté à un calendrier. Celui-ci contiendra tous les événements des fichiers importés. C’est ce qui permet de les fusionner.", $false, $false, $false, $false, $false, $true, 1, $false, "il est ajouté au calendrier final. Celui-ci contiendra tous les événements des fichiers importés. C’est ce qui permet de les fusionner. Le résultat est disponible pour l’exportation.", 2) | Out-Null

# 4/5. Remove the "Lorsque l'opération est terminée..." paragraph and the blank spacer
#      paragraph right after it, so "Améliorations" directly follows the paragraph
#      above and directly precedes the (reworded) "Le programme a atteint..." paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Lorsque l")) {
        $txt = $d.Paragraphs($i).Range.Text
        if ($txt -like "*l’opération est terminée*") {
            $target = $i
            break
        }
    }
}
if ($target -ne $null) {
    $d.Paragraphs($target + 1).Range.Delete() | Out-Null
    $d.Paragraphs($target).Range.Delete() | Out-Null
}

# 6. Reword "Le programme a atteint ses objectifs, malgré tout il serait possible d'y
#    améliorer certains points :" -> "...apporter quelques améliorations :"
$d.Content.Find.Execute("Le programme a atteint ses objectifs, malgré tout il serait possible d’y améliorer certains points", $false, $false, $false, $false, $false, $true, 1, $false, "Le programme a atteint ses objectifs, malgré tout il serait possible d’y apporter quelques améliorations", 2) | Out-Null

# 7. Merge ", ainsi il" + bookmark "_GoBack" + " sera plus facile..." into a single run,
#    dropping the _GoBack bookmark in the process.
$d.Content.Find.Execute(", ainsi il sera plus facile de fusionner de très gros calendriers. (Exemple", $false, $false, $false, $false, $false, $true, 1, $false, ", ainsi il sera plus facile de fusionner de très gros calendriers. (Exemple", 2) | Out-Null
